$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextDate($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $fmt = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = $fmt
}

# Row 10
$ws.Range("B10").Value = 4794
Set-TextDate "C10" "01/04/2020"
$ws.Range("D10").Value = 0.625
$ws.Range("E10").Value = 0.68194444444444446
$ws.Range("G10").Value = "Help debug arithmetic unit"

# Row 11
Set-TextDate "C11" "02/04/2020"
$ws.Range("D11").Value = 0.72916666666666663
$ws.Range("E11").Value = 0.77083333333333337
$ws.Range("G11").Value = "Screenshots of functional simulation waves"

# Row 12
Set-TextDate "C12" "02/04/2020"
$ws.Range("D12").Value = 0.81944444444444453
$ws.Range("E12").Value = 0.94791666666666663
$ws.Range("G12").Value = "Working on documentations and project report"

# Update view: top left cell and selection
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("G12").Select()
